$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data (row 2), pushing existing data down
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# Remove the now-trailing old rows (the window shifted, dropping the final 3 original rows)
$ws.Range("A22:C24").EntireRow.Delete()

# Populate the two newly inserted rows with the new accelerometer samples
$ws.Range("A2").Value2 = -2.025566756725311
$ws.Range("B2").Value2 = 3.52062651515007
$ws.Range("C2").Value2 = 2.27691987156868

$ws.Range("A3").Value2 = -2.230706214904786
$ws.Range("B3").Value2 = 3.561713695526123
$ws.Range("C3").Value2 = 2.031704187393189
